# Update the heading date (Find/Replace, format-preserving).
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-05-14 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-05-15 Thursday", 2)

# The practice table keeps its 20x5 shape, but every problem shifts by one
# cell (a row was removed after row 1 and a new row was inserted after the
# new row 1), so every cell's expression changes. Rewrite each cell in
# place so the existing run formatting (font/size) is preserved.
$newValues = @(
    @("2+39=", "81-14=", "61+29=", "26+31=", "15-7="),
    @("62-5=", "70-14=", "62+12=", "84+10=", "31-18="),
    @("25-23=", "6+83=", "88-69=", "53+11=", "69+12="),
    @("7+86=", "59+4=", "40+39=", "76-28=", "27+27="),
    @("5+78=", "10+55=", "60+29=", "73-69=", "41-30="),
    @("40+53=", "5+3=", "67-48=", "4+81=", "93-21="),
    @("77+6=", "31-23=", "18+26=", "70-54=", "18+60="),
    @("82-62=", "56-46=", "90-66=", "58+30=", "67+17="),
    @("71+22=", "18+39=", "38+4=", "64-29=", "36-28="),
    @("82-3=", "43-19=", "59+20=", "91+3=", "35+5="),
    @("29+20=", "71-21=", "56+4=", "99-31=", "67-49="),
    @("37+6=", "44+42=", "40+29=", "16-3=", "95-7="),
    @("95-72=", "76-33=", "66-32=", "84-13=", "78-29="),
    @("19-6=", "72-27=", "64-12=", "54+24=", "50-41="),
    @("15-0=", "5+83=", "65-34=", "1+43=", "1+62="),
    @("90-52=", "28-5=", "15+53=", "21+21=", "44-30="),
    @("49+1=", "65-28=", "28+16=", "87-44=", "45+39="),
    @("98-28=", "86-48=", "71-34=", "84-57=", "1+52="),
    @("33+4=", "1+57=", "23+39=", "86-58=", "83-38="),
    @("97-8=", "11+45=", "81-58=", "43+49=", "71+6=")
)

$t = $d.Tables.Item(1)
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$r - 1][$c - 1]
    }
}
